$d = $word.ActiveDocument

# --- Remove the last two sample paragraphs (underline-formatted and
# strike-formatted "Sample Text ..." paragraphs), keeping the subscript
# paragraph as the final one in the body.
$count = $d.Paragraphs.Count
$pUnderline = $d.Paragraphs.Item($count - 1)
$pStrike = $d.Paragraphs.Item($count)
$r = $d.Range($pUnderline.Range.Start, $pStrike.Range.End)
$r.Delete()

# --- Update eastAsia font ("DejaVu Sans" -> "Tahoma") on the styles that
# reference it: the "Normal" and "Heading" paragraph styles.
$normal = $d.Styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $d.Styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# --- Add an explicit complex-script font (w:cs="DejaVu Sans") to the
# "List", "Caption" and "Index" paragraph styles.
$list = $d.Styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $d.Styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $d.Styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
